# The document contains three occurrences of a split "<id>...</id>" marker
# where the inner id text (e.g. "p090r_a1") sits in its own run, sandwiched
# between a "<id>" run and a "</id>" run (all three in the same Courier-New
# formatting except the middle, plain-black run). The newly downloaded
# transcription replaces each with a single run containing the full,
# renumbered id tag: "<id>p090r_1</id>", "<id>p090r_2</id>", "<id>p090r_3</id>".
#
# Using Find & Replace across the whole "<id>...</id>" span merges the three
# runs into one run that takes on the formatting of the first matched
# character (the Courier-New "<id>" run), which is exactly the formatting
# the final merged run should have.

$d = $word.ActiveDocument

$pairs = @(
    @{ Old = "<id>p090r_a1</id>"; New = "<id>p090r_1</id>" },
    @{ Old = "<id>p090r_a2</id>"; New = "<id>p090r_2</id>" },
    @{ Old = "<id>p090r_a3</id>"; New = "<id>p090r_3</id>" }
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $found = $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                                  $true, 1, $false, $pair.New, 2)
    if (-not $found) {
        throw "Could not find expected text '$($pair.Old)' to replace."
    }
}
